$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New StatQuery text that replaces the old "all_studies" style query.
# It is written to the StatQuery column (C) for every data row (2-4),
# which causes the old shared string to become orphaned and the engine
# drops it, while this new text is appended as a new shared string.
$newQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Black and Tan Coonhound']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Move the active selection from B2 to B1.
$ws.Range("B1").Select()
